# Apply the emergencias.xlsx update:
#  - meta sheet: update hora_corte (B2) from 11:00 AM to 2:00 PM, move selection to B3
#  - events sheet: add new emergency row 5 (id=4, Ucayali / Padre Abad flood event),
#    move selection to H7
# Final active tab must remain "meta" (matches source sheetView tabSelected flag).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("meta")
$events = $wb.Worksheets.Item("events")

# --- events sheet: append new emergency row ---
# Row 5 previously had leftover formatting (E5:I5) from earlier rows; clear it so
# the new values land with default (unstyled) cells, like rows 2-4.
$events.Range("E5:I5").ClearFormats()

$events.Range("A5").Value = 4
$events.Range("B5").Value = "Ucayali"
$events.Range("C5").Value = "Padre Abad"
$events.Range("D5").Value = "Padre Abad"
$events.Range("E5").Value = 45988
$events.Range("F5").Value = "Lluvias intensas"
$events.Range("G5").Value = "En monitoreo"
$events.Range("H5").Value = "1 CIAI afectado (28 niños SCD)"
$events.Range("I5").Value = "Activa"
$events.Range("J5").Value = "Evaluación en curso"
$events.Range("K5").Value = -8.8943949999999994
$events.Range("L5").Value = -75.641589999999994

# --- meta sheet: hora_corte updated from 11:00 AM to 2:00 PM ---
$meta.Range("B2").Value = 0.58333333333333337

# Update selections last (select events first, then meta last so "meta" ends
# up as the active/visible tab again, matching the original workbook).
$events.Range("H7").Select()
$meta.Range("B3").Select()
